$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.61"
$ws.Range("E2").Value = "'-4.75%"
$ws.Range("G2").Value = "'2"
$ws.Range("D3").Value = "'35.02"
$ws.Range("E3").Value = "'-4.10%"
$ws.Range("G3").Value = "'2"
$ws.Range("D4").Value = "'4.925"
$ws.Range("E4").Value = "'-3.19%"
$ws.Range("G4").Value = "'2"
$ws.Range("D5").Value = "'0.07749"
$ws.Range("E5").Value = "'-3.71%"
$ws.Range("G5").Value = "'2"
$ws.Range("D6").Value = "'1.898"
$ws.Range("E6").Value = "'-12.78%"
$ws.Range("G6").Value = "'2"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.012"
$ws.Range("E7").Value = "'-3.28%"
$ws.Range("G7").Value = "'2"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'7.703"
$ws.Range("E8").Value = "'-4.04%"
$ws.Range("G8").Value = "'2"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.867"
$ws.Range("E9").Value = "'2.46%"
$ws.Range("G9").Value = "'2"
$ws.Range("B10").Value = "MXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D10").Value = "'0.9204"
$ws.Range("E10").Value = "'-0.84%"
$ws.Range("G10").Value = "'2"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.1067"
$ws.Range("E11").Value = "'6.57%"
$ws.Range("G11").Value = "'2"
$ws.Range("B12").Value = "WazirX"
$ws.Range("C12").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D12").Value = "'0.1795"
$ws.Range("E12").Value = "'-4.28%"
$ws.Range("G12").Value = "'2"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "'0.09095"
$ws.Range("E13").Value = "'-1.30%"
$ws.Range("G13").Value = "'2"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03553"
$ws.Range("E14").Value = "'-0.62%"
$ws.Range("G14").Value = "'2"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09868"
$ws.Range("E15").Value = "'-0.63%"
$ws.Range("G15").Value = "'2"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001392"
$ws.Range("E16").Value = "'-2.86%"
$ws.Range("G16").Value = "'2"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005867"
$ws.Range("E17").Value = "'3.73%"
$ws.Range("G17").Value = "'2"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.475"
$ws.Range("E18").Value = "'0.33%"
$ws.Range("G18").Value = "'2"
$ws.Range("D19").Value = "'0.3437"
$ws.Range("E19").Value = "'1.88%"
$ws.Range("G19").Value = "'2"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'5.086"
$ws.Range("E20").Value = "'0.64%"
$ws.Range("G20").Value = "'2"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1274"
$ws.Range("E21").Value = "'-7.42%"
$ws.Range("G21").Value = "'2"
$ws.Range("D22").Value = "'0.2192"
$ws.Range("G22").Value = "'2"
$ws.Range("D23").Value = "'0.04515"
$ws.Range("E23").Value = "'-1.81%"
$ws.Range("G23").Value = "'2"
$ws.Range("D24").Value = "'0.001214"
$ws.Range("E24").Value = "'-1.98%"
$ws.Range("G24").Value = "'2"
$ws.Range("D25").Value = "'0.004585"
$ws.Range("E25").Value = "'-3.67%"
$ws.Range("G25").Value = "'2"
$ws.Range("D26").Value = "'0.0001253"
$ws.Range("E26").Value = "'-3.49%"
$ws.Range("G26").Value = "'2"
$ws.Range("D27").Value = "'0.0004199"
$ws.Range("E27").Value = "'-6.54%"
$ws.Range("G27").Value = "'2"
$ws.Range("G28").Value = "'2"
$ws.Range("G29").Value = "'2"
$ws.Range("G30").Value = "'2"
$ws.Range("G31").Value = "'2"
$ws.Range("G32").Value = "'2"
$ws.Range("G33").Value = "'2"
$ws.Range("G34").Value = "'2"
$ws.Range("G35").Value = "'2"
$ws.Range("G36").Value = "'2"
$ws.Range("G37").Value = "'2"
$ws.Range("G38").Value = "'2"
$ws.Range("D39").Value = "'0.01856"
$ws.Range("E39").Value = "'-4.31%"
$ws.Range("G39").Value = "'2"
$ws.Range("D40").Value = "'0.04651"
$ws.Range("E40").Value = "'-5.95%"
$ws.Range("G40").Value = "'2"
$ws.Range("D41").Value = "'0.007578"
$ws.Range("E41").Value = "'-2.87%"
$ws.Range("G41").Value = "'2"
$ws.Range("D42").Value = "'0.009373"
$ws.Range("E42").Value = "'20.35%"
$ws.Range("G42").Value = "'2"
$ws.Range("D43").Value = "'0.1315"
$ws.Range("E43").Value = "'-6.01%"
$ws.Range("G43").Value = "'2"
$ws.Range("D44").Value = "'0.002119"
$ws.Range("E44").Value = "'1.22%"
$ws.Range("G44").Value = "'2"
$ws.Range("D45").Value = "'0.01103"
$ws.Range("E45").Value = "'-3.64%"
$ws.Range("G45").Value = "'2"
$ws.Range("D46").Value = "'0.00006015"
$ws.Range("E46").Value = "'-3.64%"
$ws.Range("G46").Value = "'2"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'0.37%"
$ws.Range("G47").Value = "'2"
$ws.Range("E48").Value = "'122.74%"
$ws.Range("G48").Value = "'2"
$ws.Range("D49").Value = "'0.001304"
$ws.Range("E49").Value = "'-31.25%"
$ws.Range("G49").Value = "'2"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.37%"
$ws.Range("G50").Value = "'2"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.37%"
$ws.Range("G51").Value = "'2"
